# Automated BRVM update (recommandations.xlsx)
# Rewrites the "Recommandations" and "Top_YTD" sheets' data rows with the
# latest recommendation figures and a newly-added row (NEI-CEDA CI moved to
# row 33, list grew by one entry).

$wb = $excel.ActiveWorkbook

$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# Each inner array: (row, Titre, JoursHausse, JoursBaisse, VariationTotale, DerniereVariation, Recommandation, Strategie)
$sheet1Rows = @(
    ,@(2, "BRVM - CONSOMMATION DE BASE     (**)", 0, 3, 737.67, 256.25, "🟡 Observer", "➖ Neutre")
    ,@(3, "BRVM-PRINCIPAL     (**)", 0, 3, 733.76, 250.62, "🟡 Observer", "➖ Neutre")
    ,@(4, "BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 3, 563.92, 193.84, "🟡 Observer", "➖ Neutre")
    ,@(5, "BRVM - INDUSTRIELS", 0, 3, 549.53, 197.02, "🟡 Observer", "➖ Neutre")
    ,@(6, "BRVM - SERVICES FINANCIERS", 0, 3, 473.34, 158.93, "🟡 Observer", "➖ Neutre")
    ,@(7, "BRVM-PRESTIGE", 0, 3, 458.48, 154.02, "🟡 Observer", "➖ Neutre")
    ,@(8, "BRVM – COMPOSITE TOTAL RETURN     (**)", 0, 3, 431.74, 145.97, "🟡 Observer", "➖ Neutre")
    ,@(9, "BRVM - ENERGIE", 0, 3, 368.62, 123.98, "🟡 Observer", "➖ Neutre")
    ,@(10, "BRVM - SERVICES PUBLICS", 0, 3, 355.24, 119.67, "🟡 Observer", "➖ Neutre")
    ,@(11, "BRVM - TELECOMMUNICATIONS", 0, 3, 296.96, 99.61, "🟡 Observer", "➖ Neutre")
    ,@(12, "EVIOSYS PACKAGING SIEM CI (SEMC)", 3, 0, 22.04, 7.43, "🟢 Achat", "✅ Renforcer")
    ,@(13, "UNIWAX CI (UNXC)", 3, 0, 21.93, 7.32, "🟢 Achat", "✅ Renforcer")
    ,@(14, "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)", 2, 0, 14.9, 7.46, "🟡 Observer", "➖ Neutre")
    ,@(15, "SETAO CI (STAC)", 2, 0, 14.48, 7.23, "🟡 Observer", "➖ Neutre")
    ,@(16, "SICABLE CI (CABC)", 1, 0, 7.5, 7.5, "🟡 Observer", "➖ Neutre")
    ,@(17, "UNILEVER CI (UNLC)", 1, 0, 7.5, 7.5, "🟡 Observer", "➖ Neutre")
    ,@(18, "SICOR CI (SICC)", 1, 0, 7.47, 7.47, "🟡 Observer", "➖ Neutre")
    ,@(19, "ORAGROUP TOGO (ORGT)", 1, 1, 4.84, 7.45, "🟡 Observer", "👀 À surveiller")
    ,@(20, "SODE CI (SDCC)", 1, 1, 0.5, -6.89, "🟡 Observer", "👀 À surveiller")
    ,@(21, "LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -1.02, -1.02, "🟡 Observer", "➖ Neutre")
    ,@(22, "BANK OF AFRICA SENEGAL (BOAS)", 0, 1, -1.59, -1.59, "🟡 Observer", "➖ Neutre")
    ,@(23, "TOTALENERGIES MARKETING CI (TTLC)", 0, 1, -1.64, -1.64, "🟡 Observer", "➖ Neutre")
    ,@(24, "BANK OF AFRICA NG (BOAN)", 0, 1, -1.7, -1.7, "🟡 Observer", "➖ Neutre")
    ,@(25, "BANK OF AFRICA BF (BOABF)", 0, 1, -1.9, -1.9, "🟡 Observer", "➖ Neutre")
    ,@(26, "CORIS BANK INTERNATIONAL (CBIBF)", 0, 1, -2, -2, "🟡 Observer", "➖ Neutre")
    ,@(27, "ONATEL BF (ONTBF)", 0, 1, -2.21, -2.21, "🟡 Observer", "➖ Neutre")
    ,@(28, "BANK OF AFRICA ML (BOAM)", 0, 1, -2.41, -2.41, "🟡 Observer", "➖ Neutre")
    ,@(29, "VIVO ENERGY CI (SHEC)", 0, 1, -2.78, -2.78, "🟡 Observer", "➖ Neutre")
    ,@(30, "ECOBANK COTE D''IVOIRE (ECOC)", 0, 1, -2.86, -2.86, "🟡 Observer", "➖ Neutre")
    ,@(31, "BERNABE CI (BNBC)", 0, 1, -3.25, -3.25, "🟡 Observer", "➖ Neutre")
    ,@(32, "SOGB CI (SOGC)", 0, 1, -4.34, -4.34, "🟡 Observer", "➖ Neutre")
    ,@(33, "NEI-CEDA CI (NEIC)", 0, 1, -7.49, -7.49, "🟡 Observer", "➖ Neutre")
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $wsReco.Cells.Item($r, 1).Value = $row[1]
    $wsReco.Cells.Item($r, 2).Value = $row[2]
    $wsReco.Cells.Item($r, 3).Value = $row[3]
    $wsReco.Cells.Item($r, 4).Value = $row[4]
    $wsReco.Cells.Item($r, 5).Value = $row[5]
    $wsReco.Cells.Item($r, 6).Value = $row[6]
    $wsReco.Cells.Item($r, 7).Value = $row[7]
}

# Each inner array: (row, Titre, Progression YTD %)
$sheet2Rows = @(
    ,@(2, "BRVM - CONSOMMATION DE BASE     (**)", 4035.18)
    ,@(3, "BRVM-PRINCIPAL     (**)", 3990.59)
    ,@(4, "BRVM - CONSOMMATION DISCRETIONNAIRE", 2287.03)
    ,@(5, "BRVM - INDUSTRIELS", 2166.04)
    ,@(6, "BRVM - SERVICES FINANCIERS", 1612.94)
    ,@(7, "BRVM-PRESTIGE", 1516.07)
    ,@(8, "BRVM – COMPOSITE TOTAL RETURN     (**)", 1351.05)
    ,@(9, "BRVM - ENERGIE", 1007.04)
    ,@(10, "BRVM - SERVICES PUBLICS", 941.9)
    ,@(11, "BRVM - TELECOMMUNICATIONS", 687.89)
)

foreach ($row in $sheet2Rows) {
    $r = $row[0]
    $wsYtd.Cells.Item($r, 1).Value = $row[1]
    $wsYtd.Cells.Item($r, 2).Value = $row[2]
}
